{"js": "// Remove the last two text paragraphs (\"Aenean justo risus...\" and\n// \"Vestibulum egestas quam massa...\") from the document body, leaving the\n// trailing empty paragraph (with the _GoBack bookmark) untouched.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  \"Aenean justo risus, hendrerit aliquam orci vel, fermentum dignissim dolor.\",\n  \"Vestibulum egestas quam massa. Sed a congue est.\"\n];\n\nfor (const p of paragraphs.items) {\n  const txt = p.text || \"\";\n  if (targets.some((t) => txt.indexOf(t) === 0)) {\n    p.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the last two text paragraphs (\"Aenean justo risus...\" and\n# \"Vestibulum egestas quam massa...\") from the document body, leaving the\n# trailing empty paragraph (with the _GoBack bookmark) untouched.\n$d = $word.ActiveDocument\n\n$targets = @(\n    \"Aenean justo risus, hendrerit aliquam orci vel, fermentum dignissim dolor.\",\n    \"Vestibulum egestas quam massa. Sed a congue est.\"\n)\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text\n    foreach ($t in $targets) {\n        if ($text.StartsWith($t)) {\n            $p.Range.Delete()\n            break\n        }\n    }\n}\n"}
